# Generate Report for handoff
# The file "7783ce40-bf73-4bd0-b0ed-5c0dede8461b.md" moves from
# "Handed back: in sync with en-US" to "Ready for handoff" status, and the
# per-locale sheets record a new "Latest Handoff Datetime" for that row.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "2016-02-16 13:56:20"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "2016-02-16 13:56:34"
